# Generate Report for Handback
# Mark handback as complete: update status text, fill in target/handback
# file columns, and stamp the handback datetime for each locale sheet.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$targetFileName = "f7a7745b-c98c-4a30-a502-df52008765f6.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn sheet -----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = $targetFileName
$wsZh.Range("J2").Value = "f7a7745b-c98c-4a30-a502-df52008765f6.0192b93ecf21300cb9d8b51f70b52313859daee9.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 15:08:02"

$wsZh.Range("I3").Value = $targetFileName
$wsZh.Range("J3").Value = "f7a7745b-c98c-4a30-a502-df52008765f6.0192b93ecf21300cb9d8b51f70b52313859daee9.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-03 15:08:02"

$zhLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4c8df79f6ef5007aa6d5b7b7cb002c3a13dd291f/e2e/f7a7745b-c98c-4a30-a502-df52008765f6.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhLinkAddress, "", "", $targetFileName)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhLinkAddress, "", "", $targetFileName)
# Match the existing custom "HyperLink" cell style (blue/underline) used by A2/A3
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276

$wsZh.Range("C1").EntireColumn.ColumnWidth = 29.14
$wsZh.Range("I1").EntireColumn.ColumnWidth = 39.14
$wsZh.Range("J1").EntireColumn.ColumnWidth = 39.14

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = $targetFileName
$wsDe.Range("J2").Value = "f7a7745b-c98c-4a30-a502-df52008765f6.0192b93ecf21300cb9d8b51f70b52313859daee9.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 15:08:15"

$wsDe.Range("I3").Value = $targetFileName
$wsDe.Range("J3").Value = "f7a7745b-c98c-4a30-a502-df52008765f6.0192b93ecf21300cb9d8b51f70b52313859daee9.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 15:08:15"

$deLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4c8df79f6ef5007aa6d5b7b7cb002c3a13dd291f/e2e/f7a7745b-c98c-4a30-a502-df52008765f6.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deLinkAddress, "", "", $targetFileName)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deLinkAddress, "", "", $targetFileName)
# Match the existing custom "HyperLink" cell style (blue/underline) used by A2/A3
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276

$wsDe.Range("C1").EntireColumn.ColumnWidth = 29.14
$wsDe.Range("I1").EntireColumn.ColumnWidth = 39.14
$wsDe.Range("J1").EntireColumn.ColumnWidth = 39.14

# --- Overview sheet status columns auto-fit -----------------------------
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 29.14
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 29.14
